$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 4338.3335
$ws.Range("I18").Value = 5385.143
$ws.Range("K18").Value = 5385.143
$ws.Range("M18").Value = -5101.143
$ws.Range("H116").Value = 25005372
$ws.Range("I116").Value = 83337830
$ws.Range("J116").Value = 5744.5713
$ws.Range("K116").Value = 83337830
$ws.Range("L116").Value = 5744.5713
$ws.Range("M116").Value = -83334388
$ws.Range("N116").Value = -12628.5713

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1528867.6
$ws.Range("I32").Value = 1648578.5
$ws.Range("K32").Value = 1648578.5
$ws.Range("M32").Value = -1648291.5
$ws.Range("H45").Value = 5732.9414
$ws.Range("I45").Value = 1714.3334
$ws.Range("J45").Value = 15377.6
$ws.Range("K45").Value = 1714.3334
$ws.Range("L45").Value = 15377.6
$ws.Range("M45").Value = -1337.3334
$ws.Range("N45").Value = -16131.6
$ws.Range("H61").Value = 3712.9858
$ws.Range("I61").Value = 1680.8983
$ws.Range("K61").Value = 1680.8983
$ws.Range("M61").Value = -1468.8983
$ws.Range("H102").Value = 999.6667
$ws.Range("I102").Value = 999.6667
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 999.6667
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 622.3333
$ws.Range("H136").Value = 3712.9858
$ws.Range("I136").Value = 1680.8983
$ws.Range("K136").Value = 5042.6949
$ws.Range("M136").Value = -2492.6949

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 9261176
$ws.Range("I20").Value = 20835454
$ws.Range("J20").Value = 1753.2
$ws.Range("K20").Value = 20835454
$ws.Range("L20").Value = 1753.2
$ws.Range("M20").Value = -20835207
$ws.Range("N20").Value = -2247.2
$ws.Range("H99").Value = 5350252
$ws.Range("I99").Value = 2577.3076
$ws.Range("J99").Value = 22730194
$ws.Range("K99").Value = 2577.3076
$ws.Range("L99").Value = 22730194
$ws.Range("M99").Value = -1079.3076
$ws.Range("N99").Value = -22733190
$ws.Range("H105").Value = 1857.9166
$ws.Range("I105").Value = 1663.1818
$ws.Range("K105").Value = 1663.1818
$ws.Range("M105").Value = 83.81819999999993

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 264.07693
$ws.Range("I22").Value = 278.66666
$ws.Range("J22").Value = 231.25
$ws.Range("K22").Value = 278.66666
$ws.Range("L22").Value = 231.25
$ws.Range("M22").Value = 71.33334000000002
$ws.Range("N22").Value = -931.25
$ws.Range("H55").Value = 50000
$ws.Range("J55").Value = 50000
$ws.Range("L55").Value = 50000
$ws.Range("N55").Value = -50630
$ws.Range("H64").Value = 65481.668
$ws.Range("J64").Value = 65481.668
$ws.Range("L64").Value = 65481.668
$ws.Range("N64").Value = -65977.66800000001
$ws.Range("H67").Value = 65481.668
$ws.Range("J67").Value = 65481.668
$ws.Range("L67").Value = 65481.668
$ws.Range("N67").Value = -67197.66800000001
$ws.Range("H99").Value = 11455.909
$ws.Range("I99").Value = 16252.75
$ws.Range("J99").Value = 8714.857
$ws.Range("K99").Value = 16252.75
$ws.Range("L99").Value = 8714.857
$ws.Range("M99").Value = -14754.75
$ws.Range("N99").Value = -11710.857
$ws.Range("H105").Value = 5956686.5
$ws.Range("I105").Value = 10205163
$ws.Range("K105").Value = 10205163
$ws.Range("M105").Value = -10203416
$ws.Range("H126").Value = 11455.909
$ws.Range("I126").Value = 16252.75
$ws.Range("J126").Value = 8714.857
$ws.Range("K126").Value = 48758.25
$ws.Range("L126").Value = 26144.571
$ws.Range("M126").Value = -46288.25
$ws.Range("N126").Value = -31084.571

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 5000528
$ws.Range("I12").Value = 1948
$ws.Range("J12").Value = 6250173
$ws.Range("K12").Value = 5844
$ws.Range("L12").Value = 18750519
$ws.Range("M12").Value = -5671
$ws.Range("N12").Value = -18750865
$ws.Range("H38").Value = 72.125
$ws.Range("J38").Value = 81.25
$ws.Range("L38").Value = 243.75
$ws.Range("N38").Value = -937.75
$ws.Range("H62").Value = 5000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("H65").Value = 5000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("H137").Value = 119495.35
$ws.Range("I137").Value = 92483.55
$ws.Range("K137").Value = 277450.65
$ws.Range("M137").Value = -272350.65

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3521.1292
$ws.Range("I132").Value = 1572.88
$ws.Range("J132").Value = 11638.833
$ws.Range("K132").Value = 4718.64
$ws.Range("L132").Value = 34916.499
$ws.Range("M132").Value = -2188.64
$ws.Range("N132").Value = -39976.499
$ws.Range("H138").Value = 89999
$ws.Range("J138").Value = 89999
$ws.Range("L138").Value = 89999
$ws.Range("N138").Value = -100279

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1626.3684
$ws.Range("I22").Value = 875
$ws.Range("K22").Value = 875
$ws.Range("M22").Value = -580
$ws.Range("H27").Value = 1626.3684
$ws.Range("I27").Value = 875
$ws.Range("K27").Value = 875
$ws.Range("M27").Value = -768
$ws.Range("H55").Value = 33333668
$ws.Range("I55").Value = 66666750
$ws.Range("K55").Value = 66666750
$ws.Range("M55").Value = -66666577
$ws.Range("H68").Value = 2826.95
$ws.Range("I68").Value = 1967.1428
$ws.Range("J68").Value = 4833.1665
$ws.Range("K68").Value = 1967.1428
$ws.Range("L68").Value = 4833.1665
$ws.Range("M68").Value = -1218.1428
$ws.Range("N68").Value = -6331.1665
$ws.Range("H71").Value = 2826.95
$ws.Range("I71").Value = 1967.1428
$ws.Range("J71").Value = 4833.1665
$ws.Range("K71").Value = 9835.714
$ws.Range("L71").Value = 24165.8325
$ws.Range("M71").Value = -6091.714
$ws.Range("N71").Value = -31653.8325
$ws.Range("H100").Value = 5837.5
$ws.Range("I100").Value = 4266.3335
$ws.Range("J100").Value = 6780.2
$ws.Range("K100").Value = 4266.3335
$ws.Range("L100").Value = 6780.2
$ws.Range("M100").Value = -3725.3335
$ws.Range("N100").Value = -7862.2

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 98183.164
$ws.Range("I122").Value = 115998.8
$ws.Range("K122").Value = 347996.4
$ws.Range("M122").Value = -345546.4
$ws.Range("H132").Value = 9623829
$ws.Range("I132").Value = 14289486
$ws.Range("K132").Value = 42868458
$ws.Range("M132").Value = -42865928
$ws.Range("H136").Value = 18541018
$ws.Range("I136").Value = 28572574
$ws.Range("K136").Value = 85717722
$ws.Range("M136").Value = -85715172

# --- Clear removed cells ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N102").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N62").ClearContents()
$ws.Range("N65").ClearContents()
